$wb = $excel.ActiveWorkbook

# --- Capabilities sheet: pCloudy endpoint device -> us ---
$wsCap = $wb.Worksheets.Item("Capabilities")
$wsCap.Cells.Item(2, 4).Value = "'https://us.pcloudy.com"
$wsCap.Cells.Item(3, 4).Value = "'https://us.pcloudy.com"

# --- DeviceList sheet: swap Android/Samsung device for iOS/Apple devices ---
$wsDev = $wb.Worksheets.Item("DeviceList")
$wsDev.Cells.Item(1, 2).Value = "APPLE_iPhone8_iOS_14.0.1_aa631"
$wsDev.Cells.Item(1, 3).Value = "APPLE_iPhone12ProMax_iOS_14.6.0_de280"
$wsDev.Cells.Item(2, 2).Value = "'14.0.1"
$wsDev.Cells.Item(2, 3).Value = "'14.6.0"
$wsDev.Cells.Item(3, 2).Value = "'pCloudyIOS"
$wsDev.Cells.Item(3, 3).Value = "'pCloudyIOS"

# --- Selection state changes ---
[void]$wsCap.Range("D7").Select()
[void]$wsDev.Select()
[void]$wsDev.Range("B13").Select()
